{"js": "// \"Getting rid of AR\" \u2014 replace the \"AR\" shorthand with \"Recommendation\"\n// (and its ${AR} merge field with ${REC}) in the heading, and swap the\n// inline \"AR\" reference in the savings paragraph with \"recommendation\".\n\nconst body = context.document.body;\n\n// 1) Heading: \"AR ${AR}\" -> \"Recommendation ${REC}\"\nconst headingHits = body.search(\"AR ${AR}\", { matchCase: true });\nheadingHits.load(\"items\");\nawait context.sync();\n\nif (headingHits.items.length === 0) {\n  throw new Error(\"Could not find 'AR ${AR}' heading text to replace.\");\n}\nheadingHits.items[0].insertText(\"Recommendation ${REC}\", \"Replace\");\nawait context.sync();\n\n// 2) Body paragraph: \"this AR is\" -> \"this recommendation is\"\nconst bodyHits = body.search(\"this AR is\", { matchCase: true });\nbodyHits.load(\"items\");\nawait context.sync();\n\nif (bodyHits.items.length === 0) {\n  throw new Error(\"Could not find 'this AR is' text to replace.\");\n}\nbodyHits.items[0].insertText(\"this recommendation is\", \"Replace\");\nawait context.sync();\n", "ps1": "# \"Getting rid of AR\" \u2014 replace the \"AR\" shorthand with \"Recommendation\"\n# (and its ${AR} merge field with ${REC}) in the heading, and swap the\n# inline \"AR\" reference in the savings paragraph with \"recommendation\".\n\n$d = $word.ActiveDocument\n\n# 1) Heading: \"AR ${AR}\" -> \"Recommendation ${REC}\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"AR `${AR}\"\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.MatchWildcards = $false\n$found1 = $find1.Execute()\nif (-not $found1) {\n    throw \"Could not find 'AR `${AR}' heading text to replace.\"\n}\n$headingRange = $d.Range($find1.Parent.Start, $find1.Parent.End)\n$headingRange.Text = \"Recommendation `${REC}\"\n\n# 2) Body paragraph: \"this AR is\" -> \"this recommendation is\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"this AR is\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$found2 = $find2.Execute()\nif (-not $found2) {\n    throw \"Could not find 'this AR is' text to replace.\"\n}\n$bodyRange = $d.Range($find2.Parent.Start, $find2.Parent.End)\n$bodyRange.Text = \"this recommendation is\"\n"}
